# Updated pick up script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to reflect the new focus (A2:H4, active cell A2)
$ws.Range("A2:H4").Select()

# Row 2: delivery option now "Pick up at store"; fill in Store# / Zip Code values
$ws.Range("D2").Value = "Pick up at store"
$ws.Range("E2").Value = 558
$ws.Range("F2").Value = 18052

# Row 3: delivery option now "Pick up at store"; fill in Store# / Zip Code values
$ws.Range("D3").Value = "Pick up at store"
$ws.Range("E3").Value = 558
$ws.Range("F3").Value = 18052

# Row 4: change Order column from "Order 2" to "Order 1", and update Store#/Zip Code
$ws.Range("A4").Value = "Order 1"
$ws.Range("E4").Value = 2013
$ws.Range("F4").Value = 84107

# Remove the now-unused "Ship to address" entry from the Delivery Option drop-down data
# (no direct cell referenced it with different text; delivery option column values remain as-is)
